# Apply odds updates to Jogos_do_Dia_Betfair_Back_Lay_2025-12-22.xlsx
# Commit message: "Atualizando o arquivo XLSX"
# Updates numeric odds/lay values across rows 2-21 (Sheet1) to match the new snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 1.24
$ws.Range("Q2").Value = 1.22
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.22
# Row 3
$ws.Range("AK3").Value = 17
$ws.Range("N3").Value = 4.4
$ws.Range("P3").Value = 2.18
$ws.Range("Q3").Value = 1.74
$ws.Range("S3").Value = 2.84
$ws.Range("T3").Value = 1.72
$ws.Range("U3").Value = 2.2
# Row 4
$ws.Range("AN4").Value = 220
$ws.Range("F4").Value = 4.7
$ws.Range("G4").Value = 6.6
$ws.Range("H4").Value = 1.89
$ws.Range("I4").Value = 2.12
$ws.Range("J4").Value = 2.74
$ws.Range("K4").Value = 3.55
$ws.Range("V4").Value = 1.89
$ws.Range("W4").Value = 1.19
# Row 5
$ws.Range("AI5").Value = 980
$ws.Range("AJ5").Value = 320
$ws.Range("AK5").Value = 160
$ws.Range("AL5").Value = 140
$ws.Range("AM5").Value = 190
$ws.Range("AN5").Value = 230
$ws.Range("H5").Value = 1.51
$ws.Range("I5").Value = 1.59
$ws.Range("J5").Value = 4.1
$ws.Range("K5").Value = 4.7
$ws.Range("L5").Value = 1.39
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 3.7
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 1.9
$ws.Range("R5").Value = 1.35
$ws.Range("S5").Value = 3.3
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 2.68
$ws.Range("X5").Value = 16
$ws.Range("Y5").Value = 8
# Row 6
$ws.Range("AJ6").Value = 130
$ws.Range("AK6").Value = 75
$ws.Range("AL6").Value = 80
$ws.Range("AM6").Value = 130
$ws.Range("AN6").Value = 85
$ws.Range("G6").Value = 5.6
$ws.Range("H6").Value = 1.82
$ws.Range("I6").Value = 2.06
$ws.Range("J6").Value = 3.25
$ws.Range("K6").Value = 4.5
$ws.Range("L6").Value = 1.35
$ws.Range("M6").Value = 1.07
$ws.Range("N6").Value = 3.15
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 1.83
$ws.Range("Q6").Value = 1.96
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 3.15
$ws.Range("U6").Value = 1.96
$ws.Range("V6").Value = 1.94
$ws.Range("Z6").Value = 14
# Row 7
$ws.Range("AI7").Value = 150
$ws.Range("F7").Value = 1.26
$ws.Range("J7").Value = 6.2
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.43
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 2.08
$ws.Range("U7").Value = 1.92
$ws.Range("V7").Value = 1.06
$ws.Range("Y7").Value = 60
# Row 8
$ws.Range("AB8").Value = 12.5
$ws.Range("AE8").Value = 60
$ws.Range("AM8").Value = 110
$ws.Range("AN8").Value = 16.5
$ws.Range("G8").Value = 2.24
$ws.Range("I8").Value = 4.6
$ws.Range("J8").Value = 3.35
$ws.Range("L8").Value = 1.31
$ws.Range("N8").Value = 3.9
$ws.Range("O8").Value = 1.27
$ws.Range("P8").Value = 2
$ws.Range("Q8").Value = 1.79
$ws.Range("R8").Value = 1.39
$ws.Range("S8").Value = 2.72
$ws.Range("U8").Value = 2.14
$ws.Range("W8").Value = 1.82
# Row 9
$ws.Range("AI9").Value = 34
$ws.Range("AJ9").Value = 70
$ws.Range("AL9").Value = 980
$ws.Range("AM9").Value = 85
$ws.Range("AN9").Value = 38
$ws.Range("L9").Value = 1.37
$ws.Range("P9").Value = 2.06
$ws.Range("Q9").Value = 1.78
$ws.Range("R9").Value = 1.42
$ws.Range("S9").Value = 2.98
$ws.Range("T9").Value = 1.68
$ws.Range("U9").Value = 2.24
# Row 10
$ws.Range("AB10").Value = 26
$ws.Range("AE10").Value = 29
$ws.Range("F10").Value = 3.5
$ws.Range("G10").Value = 5.1
$ws.Range("H10").Value = 1.83
$ws.Range("I10").Value = 2.14
$ws.Range("L10").Value = 1.28
$ws.Range("M10").Value = 1.04
$ws.Range("Q10").Value = 1.59
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 2.46
$ws.Range("U10").Value = 1.01
$ws.Range("V10").Value = 1.87
$ws.Range("W10").Value = 1.24
# Row 11
$ws.Range("AC11").Value = 19.5
$ws.Range("AD11").Value = 60
$ws.Range("AE11").Value = 320
$ws.Range("AF11").Value = 9.800000000000001
$ws.Range("AH11").Value = 38
$ws.Range("AI11").Value = 230
$ws.Range("AJ11").Value = 11.5
$ws.Range("AL11").Value = 44
$ws.Range("AM11").Value = 240
$ws.Range("AN11").Value = 4.9
$ws.Range("F11").Value = 1.23
$ws.Range("G11").Value = 1.29
$ws.Range("H11").Value = 14
$ws.Range("J11").Value = 6.4
$ws.Range("K11").Value = 7.8
$ws.Range("L11").Value = 1.23
$ws.Range("N11").Value = 5.3
$ws.Range("P11").Value = 2.5
$ws.Range("R11").Value = 1.6
$ws.Range("S11").Value = 2.22
$ws.Range("T11").Value = 2.1
$ws.Range("U11").Value = 1.72
$ws.Range("W11").Value = 4.4
$ws.Range("Z11").Value = 190
# Row 12
$ws.Range("AA12").Value = 75
$ws.Range("AB12").Value = 16
$ws.Range("AE12").Value = 50
$ws.Range("AF12").Value = 24
$ws.Range("AI12").Value = 65
$ws.Range("H12").Value = 2.98
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 3.45
$ws.Range("O12").Value = 1.25
$ws.Range("P12").Value = 1.93
$ws.Range("Q12").Value = 1.88
$ws.Range("R12").Value = 1.3
$ws.Range("S12").Value = 2.92
$ws.Range("T12").Value = 1.71
$ws.Range("U12").Value = 2.16
$ws.Range("X12").Value = 1000
$ws.Range("Y12").Value = 19
$ws.Range("Z12").Value = 32
# Row 13
$ws.Range("F13").Value = 1.34
$ws.Range("G13").Value = 1.59
$ws.Range("K13").Value = 7.8
$ws.Range("L13").Value = 1.39
$ws.Range("N13").Value = 1.56
$ws.Range("P13").Value = 1.56
$ws.Range("Q13").Value = 2.04
$ws.Range("S13").Value = 2.04
$ws.Range("W13").Value = 2.68
# Row 14
$ws.Range("AB14").Value = 18.5
$ws.Range("AC14").Value = 12
$ws.Range("AD14").Value = 16.5
$ws.Range("AE14").Value = 34
$ws.Range("AH14").Value = 18
$ws.Range("AI14").Value = 40
$ws.Range("AJ14").Value = 40
$ws.Range("AK14").Value = 28
$ws.Range("AN14").Value = 16
$ws.Range("AO14").Value = 22
$ws.Range("F14").Value = 2.28
$ws.Range("G14").Value = 2.6
$ws.Range("H14").Value = 2.78
$ws.Range("I14").Value = 3.2
$ws.Range("J14").Value = 3.35
$ws.Range("K14").Value = 4.4
$ws.Range("N14").Value = 5.1
$ws.Range("T14").Value = 1.52
$ws.Range("U14").Value = 2.5
$ws.Range("V14").Value = 1.46
$ws.Range("W14").Value = 1.62
$ws.Range("Z14").Value = 29
# Row 15
$ws.Range("F15").Value = 9.6
$ws.Range("J15").Value = 5
$ws.Range("P15").Value = 1.9
$ws.Range("Q15").Value = 1.96
# Row 16
$ws.Range("AE16").Value = 19.5
$ws.Range("AH16").Value = 20
$ws.Range("AN16").Value = 40
$ws.Range("G16").Value = 5
$ws.Range("K16").Value = 5.1
$ws.Range("P16").Value = 2.72
$ws.Range("R16").Value = 1.7
$ws.Range("S16").Value = 2.16
$ws.Range("T16").Value = 1.54
$ws.Range("Y16").Value = 17.5
$ws.Range("Z16").Value = 17.5
# Row 17
$ws.Range("F17").Value = 2.66
$ws.Range("G17").Value = 2.98
$ws.Range("H17").Value = 2.4
$ws.Range("N17").Value = 1.01
$ws.Range("Q17").Value = 1.61
$ws.Range("Y17").Value = 21
# Row 18
$ws.Range("AB18").Value = 8.199999999999999
$ws.Range("AC18").Value = 10.5
$ws.Range("AD18").Value = 30
$ws.Range("AF18").Value = 11.5
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 30
$ws.Range("AJ18").Value = 22
$ws.Range("AK18").Value = 26
$ws.Range("AL18").Value = 60
$ws.Range("AN18").Value = 18
$ws.Range("F18").Value = 1.68
$ws.Range("H18").Value = 5
$ws.Range("J18").Value = 3.2
$ws.Range("L18").Value = 1.46
$ws.Range("M18").Value = 1.07
$ws.Range("N18").Value = 2.98
$ws.Range("O18").Value = 1.41
$ws.Range("R18").Value = 1.25
$ws.Range("S18").Value = 3.7
$ws.Range("T18").Value = 2.06
$ws.Range("U18").Value = 1.75
$ws.Range("X18").Value = 14
$ws.Range("Y18").Value = 21
$ws.Range("Z18").Value = 60
# Row 19
$ws.Range("O19").Value = 1.37
# Row 20
$ws.Range("AM20").Value = 180
$ws.Range("I20").Value = 5.4
$ws.Range("O20").Value = 1.47
$ws.Range("Q20").Value = 2.42
$ws.Range("R20").Value = 1.24
$ws.Range("Y20").Value = 14
# Row 21
$ws.Range("AA21").Value = 490
$ws.Range("AB21").Value = 8.6
$ws.Range("AN21").Value = 5.7
$ws.Range("N21").Value = 4.5
$ws.Range("O21").Value = 1.24
$ws.Range("P21").Value = 2.24
$ws.Range("Q21").Value = 1.68
$ws.Range("S21").Value = 2.78
$ws.Range("T21").Value = 2.14
$ws.Range("W21").Value = 3.75

Write-Host "Applied 253 cell updates"
